$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '29.488.43'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.851.79'
$ws.Range('E3').Value = '  -0.36%  '
Set-TextValue $ws.Range('D4') '0.9985'
Set-TextValue $ws.Range('D5') '240.85'
$ws.Range('E5').Value = '  -0.72%  '
Set-TextValue $ws.Range('D6') '0.6329'
$ws.Range('E6').Value = '  -0.12%  '
Set-TextValue $ws.Range('D7') '0.9999'
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue $ws.Range('D8') '0.07573'
$ws.Range('E8').Value = '  -0.50%  '
Set-TextValue $ws.Range('D9') '0.2963'
$ws.Range('E9').Value = '  -1.15%  '
Set-TextValue $ws.Range('D10') '24.67'
$ws.Range('E10').Value = '  -0.03%  '
Set-TextValue $ws.Range('D11') '0.07709'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '1.842.67'
$ws.Range('E12').Value = '  -1.79%  '
Set-TextValue $ws.Range('D13') '5.004'
$ws.Range('E13').Value = '  -0.58%  '
Set-TextValue $ws.Range('D14') '0.6867'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('E15').Value = '  +1.42%  '
Set-TextValue $ws.Range('D16') '83.41'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '2.096.96'
$ws.Range('E17').Value = '  -3.55%  '
Set-TextValue $ws.Range('D18') '6.158'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('D19').Value = '29.500.28'
$ws.Range('E19').Value = '  -0.63%  '
Set-TextValue $ws.Range('D20') '229.13'
$ws.Range('E20').Value = '  -2.46%  '
Set-TextValue $ws.Range('D21') '12.52'
$ws.Range('E21').Value = '  -0.80%  '
Set-TextValue $ws.Range('D22') '0.9993'
$ws.Range('E22').Value = '  -0.15%  '
Set-TextValue $ws.Range('D23') '7.539'
$ws.Range('E23').Value = '  -1.91%  '
Set-TextValue $ws.Range('D24') '0.9999'
$ws.Range('E24').Value = '  -0.06%  '
Set-TextValue $ws.Range('D25') '156.89'
$ws.Range('E25').Value = '  +0.86%  '
Set-TextValue $ws.Range('D26') '0.1400'
$ws.Range('E26').Value = '  -0.11%  '
Set-TextValue $ws.Range('D27') '8.398'
$ws.Range('E27').Value = '  -1.40%  '
Set-TextValue $ws.Range('D28') '17.70'
$ws.Range('E28').Value = '  -0.49%  '
Set-TextValue $ws.Range('D29') '1.471'
$ws.Range('E29').Value = '  -0.51%  '
Set-TextValue $ws.Range('D30') '1.275'
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('E31').Value = '  -1.88%  '
Set-TextValue $ws.Range('D32') '4.133'
$ws.Range('E32').Value = '  -0.22%  '
Set-TextValue $ws.Range('D33') '4.039'
$ws.Range('E33').Value = '  -0.14%  '
Set-TextValue $ws.Range('D34') '1.848'
$ws.Range('E34').Value = '  -2.27%  '
Set-TextValue $ws.Range('D35') '1.160'
$ws.Range('E35').Value = '  -1.21%  '
Set-TextValue $ws.Range('D36') '0.7158'
$ws.Range('E36').Value = '  -1.03%  '
Set-TextValue $ws.Range('D37') '2.588'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '1.248.93'
$ws.Range('E38').Value = '  -0.89%  '
Set-TextValue $ws.Range('D39') '0.01810'
$ws.Range('E39').Value = '  -0.12%  '
Set-TextValue $ws.Range('D40') '2.780'
$ws.Range('E40').Value = '  -1.20%  '
Set-TextValue $ws.Range('D41') '0.9119'
$ws.Range('E41').Value = '  +0.60%  '
Set-TextValue $ws.Range('D42') '6.210'
$ws.Range('E42').Value = '  +0.71%  '
Set-TextValue $ws.Range('D44') '101.65'
$ws.Range('E44').Value = '  +0.00%  '
Set-TextValue $ws.Range('D45') '66.16'
$ws.Range('E45').Value = '  -2.95%  '
Set-TextValue $ws.Range('D46') '0.00000000119'
$ws.Range('E46').Value = '  +0.19%  '
Set-TextValue $ws.Range('D47') '7.094'
$ws.Range('E47').Value = '  -3.44%  '
Set-TextValue $ws.Range('D48') '0.4034'
$ws.Range('E48').Value = '  -0.61%  '
Set-TextValue $ws.Range('D49') '9.141'
$ws.Range('E49').Value = '  -0.81%  '
Set-TextValue $ws.Range('D50') '1.694'
$ws.Range('E50').Value = '  -1.27%  '
Set-TextValue $ws.Range('D51') '0.1123'
$ws.Range('E51').Value = '  -0.05%  '
